$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2026-01-22 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-23 Friday", 2) | Out-Null

# Update the division-problem answers in the first table, cell by cell (row, column)
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1,1).Range.Text = "112÷6=18, 4"
$t.Cell(1,2).Range.Text = "741÷6=123, 3"
$t.Cell(1,3).Range.Text = "550÷3=183, 1"
$t.Cell(1,4).Range.Text = "558÷4=139, 2"
$t.Cell(1,5).Range.Text = "372÷9=41, 3"

# Row 5
$t.Cell(5,1).Range.Text = "296÷2=148, 0"
$t.Cell(5,2).Range.Text = "836÷9=92, 8"
$t.Cell(5,3).Range.Text = "657÷8=82, 1"
$t.Cell(5,4).Range.Text = "398÷5=79, 3"
$t.Cell(5,5).Range.Text = "704÷8=88, 0"

# Row 9
$t.Cell(9,1).Range.Text = "994÷8=124, 2"
$t.Cell(9,2).Range.Text = "532÷3=177, 1"
$t.Cell(9,3).Range.Text = "638÷9=70, 8"
$t.Cell(9,4).Range.Text = "452÷4=113, 0"
$t.Cell(9,5).Range.Text = "756÷8=94, 4"

# Row 13
$t.Cell(13,1).Range.Text = "465÷3=155, 0"
$t.Cell(13,2).Range.Text = "710÷6=118, 2"
$t.Cell(13,3).Range.Text = "753÷3=251, 0"
$t.Cell(13,4).Range.Text = "962÷7=137, 3"
$t.Cell(13,5).Range.Text = "954÷8=119, 2"

# Row 17
$t.Cell(17,1).Range.Text = "778÷5=155, 3"
$t.Cell(17,2).Range.Text = "213÷6=35, 3"
$t.Cell(17,3).Range.Text = "961÷3=320, 1"
$t.Cell(17,4).Range.Text = "836÷6=139, 2"
$t.Cell(17,5).Range.Text = "416÷5=83, 1"
